$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GUI to MCU")
$ws2 = $wb.Worksheets.Item("MCU to GUI")

# Insert two new rows before row 12 (old row 12 "Exit parameter configuration" becomes row 14)
$ws.Rows("12:13").Insert()

# ---- Row 12: "Set voltage for EIS [mV]" ----
$ws.Range("C12").Value = "Set voltage for EIS [mV]"
$ws.Range("E12").Value = 104
$ws.Range("F12").Value = "value (4 bytes)"
$ws.Range("G12").Value = "success bool (1 byte):" + [char]10 + "0 if cmd ID not recognized"
$ws.Range("H12").Value = "Set voltage peak for the AC sinusoidal signal used for EIS"
$ws.Rows("12").RowHeight = 30

# Fix up formatting of row 12 to match the rest of the table
$ws.Range("F11:H11").Copy()
$ws.Range("F12:H12").PasteSpecial(-4122)
$ws.Range("J14").Copy()
$ws.Range("J12").PasteSpecial(-4122)

# ---- Row 13: "Set number of frequencies for EIS (int)" ----
$ws.Range("C13").Value = "Set number of frequencies for EIS (int)"
$ws.Range("E13").Value = 107
$ws.Range("F13").Value = "value (4 bytes)" + [char]10 + "+ frequencies[] (4*value bytes)"
$ws.Range("G13").Value = "success bool (1 byte):" + [char]10 + "0 if failed"
$ws.Range("H13").Value = "Set the number of frequencies used in EIS as an integer" + [char]10 + "and the individual frequencies afterwards"
$ws.Range("I13").Value = "define an arbitrary frequency array, e.g. freq[] = {100, 500, 1000}" + [char]10 + "port.sent( 107 );" + [char]10 + "port.sent( 3 ); // set value (num. of freq.)" + [char]10 + "port.sent( 100 );" + [char]10 + "port.sent( 500 );" + [char]10 + "port.sent( 1000 );"
$ws.Rows("13").RowHeight = 120

# Fix up formatting of row 13 to match similar existing cells elsewhere in the table
$ws.Range("G9").Copy()
$ws.Range("F13").PasteSpecial(-4122)
$ws.Range("G9").Copy()
$ws.Range("G13").PasteSpecial(-4122)
$ws.Range("H7").Copy()
$ws.Range("H13").PasteSpecial(-4122)
$ws2.Range("I4").Copy()
$ws.Range("I13").PasteSpecial(-4122)
$ws.Range("J14").Copy()
$ws.Range("J13").PasteSpecial(-4122)

# Merge C:D on the new rows, matching the rest of the table
$ws.Range("C12:D12").Merge() | Out-Null
$ws.Range("C13:D13").Merge() | Out-Null

# Make the bold run inside G13 ("0 if failed")
$chars = $ws.Range("G13").Characters(25, 11)
$chars.Font.Bold = $true

# Widen column F to fit the new text
$ws.Columns("F:F").ColumnWidth = 33.7109375

# Selection / active sheet bookkeeping
$ws.Range("G14").Select() | Out-Null
$ws.Activate() | Out-Null

Write-Host "done"
